$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell A4 from "PATANJALI" to the new symbol "PATANJALI-BE"
$ws.Range("A4").Value = "PATANJALI-BE"

# Update the view: scroll back to top, and select A4
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A4").Select()
